# Milestone7/crud.table.xlsx — turn the plain process list into a CRUD matrix.
# Column A keeps the "Elementary Process" labels (some renamed / reordered),
# and columns C-K get single-letter Create/Read/Update markers against the
# VENDOR / INVOICE-HEADER / INVOICE-LINE / TYPE / BRAND / INVENTORY /
# PULL-INVENTORY-LINE / PULL-INVENTORY / EMPLOYEE record types.
# Rows 33-36 of the old sheet go away, shrinking the sheet to A1:K32.

$ws = $excel.ActiveSheet

# --- Step 1: rewrite column A (the process names), rows 3-32, top to bottom ---
$ws.Range("A3").Value = "App"
$ws.Range("A4").Value = "Login"
$ws.Range("A5").Value = "Main Menu"
$ws.Range("A6").Value = "Forms Menu"
$ws.Range("A7").Value = "Reports Menu"
$ws.Range("A8").Value = "Vendor Form"
$ws.Range("A9").Value = "Type Form"
$ws.Range("A10").Value = "Brand Form "
$ws.Range("A11").Value = "Employee Form"
$ws.Range("A12").Value = "Inventory Form"
$ws.Range("A13").Value = "Invoice Form"
$ws.Range("A14").Value = "Pull-Inventory Form"
$ws.Range("A15").Value = "Inventory Report"
$ws.Range("A16").Value = "Invoice Report"
$ws.Range("A17").Value = "Pull-Inventory Report"
$ws.Range("A18").Value = "Type Query"
$ws.Range("A19").Value = "Brand Query"
$ws.Range("A20").Value = "Inventory Query"
$ws.Range("A21").Value = "Vendor Query"
$ws.Range("A22").Value = "Invoice-Line Subform"
$ws.Range("A23").Value = "Update Invoice Query"
$ws.Range("A24").Value = "Employee Query"
$ws.Range("A25").Value = "Pull-Inventory-Line Subform"
$ws.Range("A26").Value = "Update Pull-Inventory Query"
$ws.Range("A27").Value = "Inventory Report Definition"
$ws.Range("A28").Value = "Invoice Report Query"
$ws.Range("A29").Value = "Invoice Report Definition"
$ws.Range("A30").Value = "Pull-Inventory Report"
$ws.Range("A31").Value = "Update Pull-Inventory Query"
$ws.Range("A32").Value = "Pull-Inventory Report Definition"

# --- Step 2: seed the three CRUD letters once each (R, then U, then C) ---
$ws.Range("K4").Value = "R"
$ws.Range("H15").Value = "U"
$ws.Range("C8").Value = "C"

# --- Step 3: fill every other CRUD marker cell ---
$ws.Range("F9").Value = "C"
$ws.Range("G10").Value = "C"
$ws.Range("K11").Value = "C"
$ws.Range("H12").Value = "C"
$ws.Range("D13").Value = "C"
$ws.Range("E13").Value = "C"
$ws.Range("I14").Value = "C"
$ws.Range("J14").Value = "C"
$ws.Range("D16").Value = "U"
$ws.Range("E16").Value = "U"
$ws.Range("I17").Value = "U"
$ws.Range("J17").Value = "U"
$ws.Range("F18").Value = "C"
$ws.Range("G19").Value = "C"
$ws.Range("H20").Value = "C"
$ws.Range("I20").Value = "R"
$ws.Range("J20").Value = "R"
$ws.Range("C21").Value = "C"
$ws.Range("D22").Value = "C"
$ws.Range("E22").Value = "C"
$ws.Range("D23").Value = "U"
$ws.Range("E23").Value = "U"
$ws.Range("K24").Value = "C"
$ws.Range("I25").Value = "C"
$ws.Range("J25").Value = "C"
$ws.Range("I26").Value = "U"
$ws.Range("J26").Value = "U"
$ws.Range("H27").Value = "R"
$ws.Range("D28").Value = "C"
$ws.Range("E28").Value = "C"
$ws.Range("D29").Value = "R"
$ws.Range("E29").Value = "R"
$ws.Range("I30").Value = "U"
$ws.Range("J30").Value = "U"
$ws.Range("I31").Value = "U"
$ws.Range("J31").Value = "U"
$ws.Range("I32").Value = "R"
$ws.Range("J32").Value = "R"

# --- Step 4: the old sheet had 4 extra rows (33-36) past the new bottom row (32) ---
$ws.Range("A33").Value = ""
$ws.Range("A34").Value = ""
$ws.Range("A35").Value = ""
$ws.Range("A36").Value = ""

# --- Step 5: move the selection from K1 to C5, like the saved file shows ---
$null = $ws.Range("C5").Select()

# --- Step 6: match the saved window geometry, best-effort (host may not
#     persist window pixel geometry into the saved bookViews element) ---
try {
    $excel.ActiveWindow.Width = 20490
    $excel.ActiveWindow.Height = 7530
} catch {
}
